$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row for release/8.0.17, mirroring the existing table layout
# (A: release/branch name, B-E: status marks for dev/sit/uat/pre-prod/prod)
$ws.Range("A20").Value = "release/8.0.17"
$ws.Range("B20").Value = "X"
$ws.Range("C20").Value = "X"
$ws.Range("D20").Value = "X"
$ws.Range("E20").Value = "X"

# The previous last row (19) has no explicit cell style, so match that by
# resetting the new row back to the default "Normal" style (columns carry a
# style-2 default that would otherwise bleed into the freshly written cells).
$ws.Range("A20:E20").Style = "Normal"
